$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H64").Value = 3696.037
$ws.Range("I64").Value = 3330.7856
$ws.Range("J64").Value = 4089.3845
$ws.Range("K64").Value = 3330.7856
$ws.Range("L64").Value = 4089.3845
$ws.Range("M64").Value = -3082.7856
$ws.Range("N64").Value = -4585.3845
$ws.Range("H67").Value = 3696.037
$ws.Range("I67").Value = 3330.7856
$ws.Range("J67").Value = 4089.3845
$ws.Range("K67").Value = 3330.7856
$ws.Range("L67").Value = 4089.3845
$ws.Range("M67").Value = -2472.7856
$ws.Range("N67").Value = -5805.3845
$ws.Range("H137").Value = 1689.1724
$ws.Range("I137").Value = 1663.5641
$ws.Range("K137").Value = 4990.692300000001
$ws.Range("M137").Value = -2440.692300000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2037.6666
$ws.Range("I2").Value = 1229.7368
$ws.Range("J2").Value = 5107.8
$ws.Range("K2").Value = 1229.7368
$ws.Range("L2").Value = 5107.8
$ws.Range("M2").Value = -1116.7368
$ws.Range("N2").Value = -5333.8
$ws.Range("H61").Value = 4442.125
$ws.Range("I61").Value = 3580
$ws.Range("J61").Value = 5112.6665
$ws.Range("K61").Value = 3580
$ws.Range("L61").Value = 5112.6665
$ws.Range("M61").Value = -3368
$ws.Range("N61").Value = -5536.6665
$ws.Range("H74").Value = 11737.333
$ws.Range("I74").Value = 12884.8
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 12884.8
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -12010.8
$ws.Range("N74").Value = -7748
$ws.Range("H77").Value = 11737.333
$ws.Range("I77").Value = 12884.8
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 64424
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -60056
$ws.Range("N77").Value = -38736
$ws.Range("H116").Value = 2037.6666
$ws.Range("I116").Value = 1229.7368
$ws.Range("J116").Value = 5107.8
$ws.Range("K116").Value = 1229.7368
$ws.Range("L116").Value = 5107.8
$ws.Range("M116").Value = 1064.2632
$ws.Range("N116").Value = -9695.799999999999
$ws.Range("H132").Value = 7755.2285
$ws.Range("I132").Value = 5725.875
$ws.Range("J132").Value = 12182.909
$ws.Range("K132").Value = 17177.625
$ws.Range("L132").Value = 36548.727
$ws.Range("M132").Value = -14647.625
$ws.Range("N132").Value = -41608.727
$ws.Range("H136").Value = 4442.125
$ws.Range("I136").Value = 3580
$ws.Range("J136").Value = 5112.6665
$ws.Range("K136").Value = 10740
$ws.Range("L136").Value = 15337.9995
$ws.Range("M136").Value = -8190
$ws.Range("N136").Value = -20437.9995

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2037.6666
$ws.Range("I3").Value = 1229.7368
$ws.Range("J3").Value = 5107.8
$ws.Range("K3").Value = 1229.7368
$ws.Range("L3").Value = 5107.8
$ws.Range("M3").Value = -1115.7368
$ws.Range("N3").Value = -5335.8
$ws.Range("H134").Value = 6995.614
$ws.Range("I134").Value = 6639.9546
$ws.Range("J134").Value = 7351.273
$ws.Range("K134").Value = 19919.8638
$ws.Range("L134").Value = 22053.819
$ws.Range("M134").Value = -17384.8638
$ws.Range("N134").Value = -27123.819
$ws.Range("H139").Value = 74733.336
$ws.Range("J139").Value = 74733.336
$ws.Range("L139").Value = 74733.336
$ws.Range("N139").Value = -85013.336

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 2164.3333
$ws.Range("I6").Value = 359.66666
$ws.Range("K6").Value = 1078.99998
$ws.Range("M6").Value = -965.9999800000001
$ws.Range("I12").Value = 58.333332
$ws.Range("J12").Value = 191.10527
$ws.Range("K12").Value = 174.999996
$ws.Range("L12").Value = 573.3158099999999
$ws.Range("M12").Value = -1.99999600000001
$ws.Range("N12").Value = -919.3158099999999
$ws.Range("H42").Value = 3150
$ws.Range("J42").Value = 3150
$ws.Range("L42").Value = 9450
$ws.Range("N42").Value = -10518
$ws.Range("H46").Value = 741.6667
$ws.Range("I46").Value = 250
$ws.Range("J46").Value = 987.5
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 2962.5
$ws.Range("M46").Value = -659
$ws.Range("N46").Value = -3144.5
$ws.Range("H51").Value = 3512.5
$ws.Range("I51").Value = 300
$ws.Range("K51").Value = 900
$ws.Range("M51").Value = -440
$ws.Range("H55").Value = 1786.4286
$ws.Range("J55").Value = 1984.1666
$ws.Range("L55").Value = 5952.4998
$ws.Range("N55").Value = -6306.4998
$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25372
$ws.Range("H64").Value = 3876.5
$ws.Range("I64").Value = 506
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 1518
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -1248
$ws.Range("N64").Value = -15540
$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78864
$ws.Range("H67").Value = 3876.5
$ws.Range("I67").Value = 506
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 1518
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -582
$ws.Range("N67").Value = -16872
$ws.Range("H80").Value = 3990
$ws.Range("J80").Value = 3990
$ws.Range("L80").Value = 11970
$ws.Range("N80").Value = -13842
$ws.Range("H83").Value = 3990
$ws.Range("J83").Value = 3990
$ws.Range("L83").Value = 35910
$ws.Range("N83").Value = -45270
$ws.Range("H137").Value = 2986.875
$ws.Range("I137").Value = 1982.5
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 5947.5
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -847.5
$ws.Range("N137").Value = -28200

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2860.0476
$ws.Range("I80").Value = 2610
$ws.Range("J80").Value = 3193.4443
$ws.Range("K80").Value = 2610
$ws.Range("L80").Value = 3193.4443
$ws.Range("M80").Value = -1612
$ws.Range("N80").Value = -5189.4443
$ws.Range("H83").Value = 2860.0476
$ws.Range("I83").Value = 2610
$ws.Range("J83").Value = 3193.4443
$ws.Range("K83").Value = 13050
$ws.Range("L83").Value = 15967.2215
$ws.Range("M83").Value = -8058
$ws.Range("N83").Value = -25951.2215
$ws.Range("H132").Value = 7782.5884
$ws.Range("I132").Value = 12256.889
$ws.Range("J132").Value = 2749
$ws.Range("K132").Value = 36770.667
$ws.Range("L132").Value = 8247
$ws.Range("M132").Value = -34240.667
$ws.Range("N132").Value = -13307

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1182.862
$ws.Range("I22").Value = 1203.8948
$ws.Range("J22").Value = 1142.9
$ws.Range("K22").Value = 1203.8948
$ws.Range("L22").Value = 1142.9
$ws.Range("M22").Value = -908.8948
$ws.Range("N22").Value = -1732.9
$ws.Range("H27").Value = 1182.862
$ws.Range("I27").Value = 1203.8948
$ws.Range("J27").Value = 1142.9
$ws.Range("K27").Value = 1203.8948
$ws.Range("L27").Value = 1142.9
$ws.Range("M27").Value = -1096.8948
$ws.Range("N27").Value = -1356.9
$ws.Range("H46").Value = 834636.75
$ws.Range("I46").Value = 875
$ws.Range("J46").Value = 1668398.5
$ws.Range("K46").Value = 875
$ws.Range("L46").Value = 1668398.5
$ws.Range("M46").Value = -687
$ws.Range("N46").Value = -1668774.5
$ws.Range("H68").Value = 3640
$ws.Range("I68").Value = 3250
$ws.Range("J68").Value = 3796
$ws.Range("K68").Value = 3250
$ws.Range("L68").Value = 3796
$ws.Range("M68").Value = -2501
$ws.Range("N68").Value = -5294
$ws.Range("H71").Value = 3640
$ws.Range("I71").Value = 3250
$ws.Range("J71").Value = 3796
$ws.Range("K71").Value = 16250
$ws.Range("L71").Value = 18980
$ws.Range("M71").Value = -12506
$ws.Range("N71").Value = -26468
$ws.Range("H132").Value = 76927970
$ws.Range("I132").Value = 111113180
$ws.Range("J132").Value = 11249.5
$ws.Range("K132").Value = 333339540
$ws.Range("L132").Value = 33748.5
$ws.Range("M132").Value = -333337010
$ws.Range("N132").Value = -38808.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5435687.5
$ws.Range("I136").Value = 5814712
$ws.Range("J136").Value = 2998.3333
$ws.Range("K136").Value = 17444136
$ws.Range("L136").Value = 8994.999899999999
$ws.Range("M136").Value = -17441586
$ws.Range("N136").Value = -14094.9999
